# Apply the dataset corrections described in the commit:
#  - normalize a couple of first/last names and company names
#  - fill in contact-match details (CRM lookup results) for rows 2 and 14
#    now that the lookup logic finds an existing contact instead of
#    recommending a brand-new one
#  - tweak a couple of confidence scores

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Gosia Owerko -> Malgorzata Owerko) ---
$ws.Range("E2").Value = "Malgorzata"
$ws.Range("O2").Value = "TAK"
$ws.Range("P2").Value = "'751364000041342001"
$ws.Range("Q2").Value = "Malgorzata Owerko"
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = "E+L+F"
$ws.Range("U2").Value = 1
$ws.Range("AA2").Value = "link_to_existing"
$ws.Range("AB2").Value = 1
$ws.Range("AC2").Value = "Znaleziono istniejący kontakt: Malgorzata Owerko (E+L+F)"

# --- Row 7 (Magda -> Magdalena) ---
$ws.Range("E7").Value = "Magdalena"

# --- Row 9 (confidence tweak) ---
$ws.Range("AB9").Value = 0.75

# --- Row 14 (Dabrowski -> Dąbrowski, contact match found) ---
$ws.Range("F14").Value = "Dąbrowski"
$ws.Range("O14").Value = "TAK"
$ws.Range("P14").Value = "'751364000002328390"
$ws.Range("Q14").Value = "Michał Dąbrowski"
$ws.Range("R14").Value = 4
$ws.Range("S14").Value = "E+L+F"
$ws.Range("U14").Value = 3
$ws.Range("AA14").Value = "link_to_existing"
$ws.Range("AB14").Value = 1
$ws.Range("AC14").Value = "Znaleziono istniejący kontakt: Michał Dąbrowski (E+L+F)"

# --- Row 15 (HVD Holding -> HVD) ---
$ws.Range("I15").Value = "HVD"

# --- Row 16 (confidence tweak) ---
$ws.Range("AB16").Value = 0.75
